# Finalise wet weather data and conductivity conversion
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capitalise the conductivity / specific-conductivity conversion names
$ws.Range("B3").Value = "COND"
$ws.Range("B4").Value = "spCOND"

# Leave selection on B4, matching the finalised state of the sheet
$ws.Range("B4").Select()
